$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 283-284, pushing the existing rows 283:299 down to 285:301
$ws.Range("283:284").Insert()

# Row 283 - new weekly entry (Primera)
$ws.Cells.Item(283, 1).Value = 9
$ws.Cells.Item(283, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(283, 3).Value = "Metropolitana"
$ws.Cells.Item(283, 4).Value = 44516
$ws.Cells.Item(283, 5).Value = 13
$ws.Cells.Item(283, 6).Value = 100114014
$ws.Cells.Item(283, 7).Value = "Betarraga"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 6100
$ws.Cells.Item(283, 11).Value = 90
$ws.Cells.Item(283, 12).Value = 100
$ws.Cells.Item(283, 13).Value = 95
$ws.Cells.Item(283, 14).Value = "`$/unidad"
$ws.Cells.Item(283, 15).Value = "Región Metropolitana"
$ws.Cells.Item(283, 16).Value = 95
$ws.Cells.Item(283, 17).Value = 1
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# Row 284 - new weekly entry (Segunda)
$ws.Cells.Item(284, 1).Value = 9
$ws.Cells.Item(284, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(284, 3).Value = "Metropolitana"
$ws.Cells.Item(284, 4).Value = 44516
$ws.Cells.Item(284, 5).Value = 13
$ws.Cells.Item(284, 6).Value = 100114014
$ws.Cells.Item(284, 7).Value = "Betarraga"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Segunda"
$ws.Cells.Item(284, 10).Value = 2500
$ws.Cells.Item(284, 11).Value = 60
$ws.Cells.Item(284, 12).Value = 70
$ws.Cells.Item(284, 13).Value = 65
$ws.Cells.Item(284, 14).Value = "`$/unidad"
$ws.Cells.Item(284, 15).Value = "Región Metropolitana"
$ws.Cells.Item(284, 16).Value = 65
$ws.Cells.Item(284, 17).Value = 1
$ws.Cells.Item(284, 18).Value = "Hortaliza"
